# edit.ps1 - Applies "Session 32 Cars Start" changes to Session 30.docx
$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1) "Real Information" block -> Session Time / Session Date cell:
#    add "5:00PM - 6:00PM AWST" after the leading space, and "11 May 2021"
#    after the manual line break.
# ---------------------------------------------------------------------------
$t3 = $d.Tables.Item(3)
$row5 = $t3.Rows.Item(5)
$timeCell = $row5.Cells.Item(2)
$row6 = $t3.Rows.Item(6)
$nextCellStart = $row6.Cells.Item(1).Range.Start
$timeRange = $d.Range($timeCell.Range.Start, $nextCellStart)
$timeRange.InsertXML("<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:ind w:left='0'/></w:pPr><w:r><w:t xml:space='preserve'> </w:t></w:r><w:r><w:t>5:00PM – 6:00PM AWST</w:t></w:r><w:r><w:br/></w:r><w:r><w:t>11 May 2021</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 2) "Github Name" label cell: drop the spell-check wrapper runs and merge
#    into a single run.
# ---------------------------------------------------------------------------
$t3 = $d.Tables.Item(3)
$row6 = $t3.Rows.Item(6)
$githubCell = $row6.Cells.Item(1)
$row7 = $t3.Rows.Item(7)
$afterGithubStart = $row7.Cells.Item(1).Range.Start
$githubRange = $d.Range($githubCell.Range.Start, $afterGithubStart)
$githubRange.InsertXML("<w:p $ns><w:pPr><w:pStyle w:val='ListParagraph'/><w:ind w:left='0'/></w:pPr><w:r><w:t>Github Name</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 3) "Activities Completed" checklist (nested table): update the three
#    checklist items' text, and give the third row an explicit row height.
# ---------------------------------------------------------------------------
$t3 = $d.Tables.Item(3)
$row8 = $t3.Rows.Item(8)
$checklistCell = $row8.Cells.Item(2)
$row9 = $t3.Rows.Item(9)
$afterChecklistStart = $row9.Cells.Item(1).Range.Start
$checklistRange = $d.Range($checklistCell.Range.Start, $afterChecklistStart)
$checklistRange.InsertXML('<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblBorders><w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:insideH w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:insideV w:val="none" w:sz="0" w:space="0" w:color="auto"/></w:tblBorders><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="343"/><w:gridCol w:w="5627"/></w:tblGrid><w:tr w:rsidR="00990AD9" w14:paraId="1BC04DEF" w14:textId="77777777" w:rsidTr="00053B1C"><w:tc><w:tcPr><w:tcW w:w="343" w:type="dxa"/></w:tcPr><w:p w14:paraId="5DEE0667" w14:textId="43D646C7" w:rsidR="00990AD9" w:rsidRDefault="00990AD9" w:rsidP="00604FC1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="0"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>□</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5627" w:type="dxa"/></w:tcPr><w:p w14:paraId="0F7878D5" w14:textId="0D654633" w:rsidR="00990AD9" w:rsidRDefault="006C2F8D" w:rsidP="00604FC1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="0"/></w:pPr><w:r><w:t>Created basic terrain</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C2F8D" w14:paraId="69E32608" w14:textId="77777777" w:rsidTr="00053B1C"><w:tc><w:tcPr><w:tcW w:w="343" w:type="dxa"/></w:tcPr><w:p w14:paraId="1218479B" w14:textId="434896A3" w:rsidR="006C2F8D" w:rsidRDefault="006C2F8D" w:rsidP="00604FC1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="0"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>□</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5627" w:type="dxa"/></w:tcPr><w:p w14:paraId="487FCCDD" w14:textId="1FB1C979" w:rsidR="006C2F8D" w:rsidRDefault="006C2F8D" w:rsidP="00604FC1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="0"/></w:pPr><w:r><w:t>Made grass and road textures</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="006C2F8D" w14:paraId="0C958D6C" w14:textId="77777777" w:rsidTr="00053B1C"><w:trPr><w:trHeight w:val="60"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="343" w:type="dxa"/></w:tcPr><w:p w14:paraId="79F480AF" w14:textId="52C7BCA9" w:rsidR="006C2F8D" w:rsidRDefault="006C2F8D" w:rsidP="00604FC1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="0"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>□</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5627" w:type="dxa"/></w:tcPr><w:p w14:paraId="33F9AE04" w14:textId="655CD2F4" w:rsidR="006C2F8D" w:rsidRDefault="006C2F8D" w:rsidP="00604FC1"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="0"/></w:pPr><w:r><w:t>Made mountain with winding road</w:t></w:r></w:p></w:tc></w:tr></w:tbl>')

# ---------------------------------------------------------------------------
# 4) Footer: merge the three runs (incl. the spell-checked "Openic") into a
#    single run with the full footer text.
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ftrRange = $ftr.Range
$ftrRange.InsertXML("<w:p $ns><w:pPr><w:pStyle w:val='Footer'/><w:jc w:val='center'/></w:pPr><w:r><w:t>Game Development Tutoring | Openic Development © All rights reserved 2021</w:t></w:r></w:p>")
